$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the office id (numeric value) in A2
$ws.Range("A2").Value = 803736

# Update the Chinese office title in B2
$ws.Range("B2").Value = "押宴官"

# Update the English translation in D2
$ws.Range("D2").Value = "Banquet Master of Ceremonies(Tackett)"

# Update the pinyin in E2
$ws.Range("E2").Value = "ya yan guan"
